$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.218.83"
$ws.Range("E2").Value = "  +0.96%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.476.60"
$ws.Range("E3").Value = "  +1.25%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.13"
$ws.Range("E5").Value = "  +0.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.72"
$ws.Range("E6").Value = "  +0.62%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.39%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.555"
$ws.Range("E8").Value = "  +0.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.502.66"
$ws.Range("E9").Value = "  +2.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0971"
$ws.Range("E10").Value = "  -0.48%  "

$ws.Range("E11").Value = "  -0.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.13"
$ws.Range("E12").Value = "  -1.75%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.330"
$ws.Range("E13").Value = "  -2.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.919.42"
$ws.Range("E14").Value = "  +1.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.181.38"
$ws.Range("E15").Value = "  +0.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.04"
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000135"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.487.44"
$ws.Range("E18").Value = "  +1.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.69"
$ws.Range("E19").Value = "  +1.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "320.38"
$ws.Range("E20").Value = "  +0.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.15"
$ws.Range("E21").Value = "  +0.49%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.02"
$ws.Range("E22").Value = "  +5.78%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.79"
$ws.Range("E24").Value = "  -0.38%  "

$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.399"
$ws.Range("E25").Value = "  -1.15%  "

$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.993"
$ws.Range("E26").Value = "  -0.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("E27").Value = "  +0.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.32"
$ws.Range("E28").Value = "  +0.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0750"
$ws.Range("E29").Value = "  +2.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.18"
$ws.Range("E30").Value = "  +1.26%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  +3.72%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.70"
$ws.Range("E32").Value = "  +1.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.20"
$ws.Range("E33").Value = "  +0.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.992"
$ws.Range("E35").Value = "  -0.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.99"
$ws.Range("E36").Value = "  +0.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.25"
$ws.Range("E37").Value = "  -2.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.93"
$ws.Range("E38").Value = "  +0.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.67"
$ws.Range("E39").Value = "  +1.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.45"
$ws.Range("E40").Value = "  -0.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.771"
$ws.Range("E41").Value = "  -1.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "275.70"
$ws.Range("E42").Value = "  +2.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.06"
$ws.Range("E43").Value = "  +2.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.42"
$ws.Range("E44").Value = "  +0.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.593"
$ws.Range("E45").Value = "  +1.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0918"
$ws.Range("E46").Value = "  +1.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.66"
$ws.Range("E47").Value = "  -1.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0498"
$ws.Range("E48").Value = "  +3.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.67"
$ws.Range("E49").Value = "  +1.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0212"
$ws.Range("E50").Value = "  +1.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.82"
$ws.Range("E51").Value = "  +1.80%  "
